$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Clear entire used range of sheet1 first so old layout (A1:Q5) is wiped ---
$ws1.Range("A1:Q5").Clear()

# --- New header row (no style) ---
$ws1.Cells.Item(1,1).Value  = "발주일자"
$ws1.Cells.Item(1,2).Value  = "납기일자"
$ws1.Cells.Item(1,3).Value  = "거래처명"
$ws1.Cells.Item(1,4).Value  = "거래처 이메일"
$ws1.Cells.Item(1,5).Value  = "납품처명"
$ws1.Cells.Item(1,6).Value  = "납품처 이메일"
$ws1.Cells.Item(1,7).Value  = "프로젝트명"
$ws1.Cells.Item(1,8).Value  = "대분류"
$ws1.Cells.Item(1,9).Value  = "중분류"
$ws1.Cells.Item(1,10).Value = "소분류"
$ws1.Cells.Item(1,11).Value = "품목명"
$ws1.Cells.Item(1,12).Value = "규격"
$ws1.Cells.Item(1,13).Value = "수량"
$ws1.Cells.Item(1,14).Value = "단가"
$ws1.Cells.Item(1,15).Value = "총금액"
$ws1.Cells.Item(1,16).Value = "비고"

# Force date-like text columns (A, B) to be stored as plain text, not dates
$ws1.Range("A2:B5").NumberFormat = "@"

# --- Row 2 ---
$ws1.Cells.Item(2,1).Value  = "2025-08-30"
$ws1.Cells.Item(2,2).Value  = "2025-09-07"
$ws1.Cells.Item(2,3).Value  = "티에스이앤씨"
$ws1.Cells.Item(2,4).Value  = "티에스이앤씨@example.com"
$ws1.Cells.Item(2,5).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2,6).Value  = "delivery@example.com"
$ws1.Cells.Item(2,7).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(2,8).Value  = "6. 안전관리비"
$ws1.Cells.Item(2,9).Value  = "1) 안전장비"
$ws1.Cells.Item(2,10).Value = "기타"
$ws1.Cells.Item(2,11).Value = "신호봉(건전지포함)"
$ws1.Cells.Item(2,12).Value = "KS규격-1"
$ws1.Cells.Item(2,13).Value = 2
$ws1.Cells.Item(2,14).Value = 5000
$ws1.Cells.Item(2,15).Value = 11000

# --- Row 3 ---
$ws1.Cells.Item(3,1).Value  = "2025-09-10"
$ws1.Cells.Item(3,2).Value  = "2025-09-24"
$ws1.Cells.Item(3,3).Value  = "티에스이앤씨"
$ws1.Cells.Item(3,4).Value  = "티에스이앤씨@example.com"
$ws1.Cells.Item(3,5).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3,6).Value  = "delivery@example.com"
$ws1.Cells.Item(3,7).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(3,8).Value  = "6. 안전관리비"
$ws1.Cells.Item(3,9).Value  = "1) 안전장비"
$ws1.Cells.Item(3,10).Value = "기타"
$ws1.Cells.Item(3,11).Value = "안전 1차 - 탄화포"
$ws1.Cells.Item(3,12).Value = "KS규격-2"
$ws1.Cells.Item(3,13).Value = 1
$ws1.Cells.Item(3,14).Value = 120000
$ws1.Cells.Item(3,15).Value = 132000
$ws1.Cells.Item(3,16).Value = "1차"

# --- Row 4 ---
$ws1.Cells.Item(4,1).Value  = "2025-09-15"
$ws1.Cells.Item(4,2).Value  = "2025-09-09"
$ws1.Cells.Item(4,3).Value  = "티에스이앤씨"
$ws1.Cells.Item(4,4).Value  = "티에스이앤씨@example.com"
$ws1.Cells.Item(4,5).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4,6).Value  = "delivery@example.com"
$ws1.Cells.Item(4,7).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(4,8).Value  = "2. 부자재비"
$ws1.Cells.Item(4,9).Value  = "2) 창호"
$ws1.Cells.Item(4,10).Value = "기타"
$ws1.Cells.Item(4,11).Value = "고임목 12T"
$ws1.Cells.Item(4,12).Value = "KS규격-3"
$ws1.Cells.Item(4,13).Value = 1
$ws1.Cells.Item(4,14).Value = 50000
$ws1.Cells.Item(4,15).Value = 55000

# --- Row 5 ---
$ws1.Cells.Item(5,1).Value  = "2025-08-22"
$ws1.Cells.Item(5,2).Value  = "2025-09-27"
$ws1.Cells.Item(5,3).Value  = "티에스이앤씨"
$ws1.Cells.Item(5,4).Value  = "티에스이앤씨@example.com"
$ws1.Cells.Item(5,5).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(5,6).Value  = "delivery@example.com"
$ws1.Cells.Item(5,7).Value  = "힐스테이트 도곡동1차"
$ws1.Cells.Item(5,8).Value  = "2. 부자재비"
$ws1.Cells.Item(5,9).Value  = "2) 창호"
$ws1.Cells.Item(5,10).Value = "기타"
$ws1.Cells.Item(5,11).Value = "고임목 혼합"
$ws1.Cells.Item(5,12).Value = "KS규격-4"
$ws1.Cells.Item(5,13).Value = 2
$ws1.Cells.Item(5,14).Value = 50000
$ws1.Cells.Item(5,15).Value = 110000

# Drop the temporary text-number-format now that the literal strings are locked in,
# so the date-like cells end up with no explicit style (matches original formatting).
$ws1.Range("A2:B5").ClearFormats()

# --- Fix sheet2 (갑지) and sheet3 (을지): remove empty remark cells I2, I4, I5 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,9).ClearContents()
$ws2.Cells.Item(4,9).ClearContents()
$ws2.Cells.Item(5,9).ClearContents()

$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,9).ClearContents()
$ws3.Cells.Item(4,9).ClearContents()
$ws3.Cells.Item(5,9).ClearContents()

Write-Host "done"
